# Reviewer pass #1 edits: rename/clean up the header labels in column A
# (and the mirrored label in AR1), then leave the selection on AR2,
# matching the order in which the author touched the cells so that the
# shared-string table is rebuilt in the same sequence.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value  = "podyplomowe"
$ws.Range("A6").Value  = "stan cywilny"
$ws.Range("A8").Value  = "Q1 staż"
$ws.Range("A2").Value  = "płeć"
$ws.Range("A4").Value  = "szkoła"
$ws.Range("A9").Value  = "Q2 etaty"
$ws.Range("A10").Value = "Q3 godz"
$ws.Range("A11").Value = "Q4 system"
$ws.Range("A12").Value = "Q5 zadowolona"
$ws.Range("A13").Value = "Q6 staż na 1 oddz."
$ws.Range("A14").Value = "Q7 wspierajace"
$ws.Range("A16").Value = "Q9 syst. motyw."
$ws.Range("A17").Value = "Q10 do domu emocje"
$ws.Range("A18").Value = "Q11 traumatyzacja"
$ws.Range("A19").Value = "Q12 potrtafi niwelowac"
$ws.Range("A20").Value = "Q13 używki"
$ws.Range("A21").Value = "Q14 dylemat rozwoj"
$ws.Range("A22").Value = "Q15 presja społeczna"
$ws.Range("A23").Value = "Q16 doswiad. przemocy"
$ws.Range("A24").Value = "Q17 post. dochodz."
$ws.Range("A25").Value = "Q18 bezp. covid"
$ws.Range("A26").Value = "Q19 popiera strajki"
$ws.Range("A27").Value = "Q20 satys. autonomi"
$ws.Range("A28").Value = "Q21 kompetencje"
$ws.Range("A29").Value = "Q22 potrz. psychologa"
$ws.Range("A30").Value = "Q23 posiada pasje"
$ws.Range("A31").Value = "Q24 uczest w życiu"
$ws.Range("A32").Value = "Q25 partner w SOZ"
$ws.Range("A33").Value = "Q26 rodzina wspiera"
$ws.Range("A34").Value = "Q27 konflikt"
$ws.Range("A35").Value = "Q28 praca pow niezgody"
$ws.Range("A36").Value = "Q29 rozpad związku"
$ws.Range("A37").Value = "Q30 rozważ. rezygnację"
$ws.Range("A38").Value = "Q31 pracuje mimo ch."
$ws.Range("A39").Value = "Q32 korzysta ZLA"
$ws.Range("A40").Value = "Q33 choruje"
$ws.Range("A41").Value = "Q34 empatia"
$ws.Range("A42").Value = "Q35 wypal. zawod."
$ws.Range("A43").Value = "Q36 dobry wybór"
$ws.Range("A44").Value = "Q37 dziecko piel."
$ws.Range("AR1").Value = "Q37"

$ws.Range("AR2").Select()
